$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (45203 -> 45204) for every data row (rows 2 through 537).
$ws.Range("C2:C537").Value = 45204
